$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated sval data (filtered save games) - literal values per diff
$values = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732; G = 6.201049113329182 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    4 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 5.553084769722144 }
    6 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732; G = 3.781711156805759 }
    7 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732; G = 9.295990156953671 }
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("B$row").Value = $v.B
    $ws.Range("C$row").Value = $v.C
    $ws.Range("D$row").Value = $v.D
    $ws.Range("E$row").Value = $v.E
    $ws.Range("G$row").Value = $v.G
}
